$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "correct" emotion label (sadness) for rows 20 and 21
$ws.Range("C20").Value = "sadness"
$ws.Range("C21").Value = "sadness"

# Update the active selection to match the recorded edit position
$ws.Range("D21").Select()
